$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) About sheet: fix typo "in the Brazil" -> "in Brazil", and push the
#    Notes block (rows 4-6) down to rows 8-10 (4 blank rows inserted).
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$noteVal   = $wsAbout.Range("A4").Value()
$ldvVal    = $wsAbout.Range("A5").Value()
$brazilVal = "in Brazil for any vehicle type."

$wsAbout.Rows("4:6").Delete()

$wsAbout.Range("A8").Value = $noteVal
$wsAbout.Range("A8").Font.Bold = $true
$wsAbout.Range("A9").Value = $ldvVal
$wsAbout.Range("A10").Value = $brazilVal

[void]$wsAbout.Range("A11").Select()

# ---------------------------------------------------------------------
# 2) BMRESP-passenger / BMRESP-freight: first forecast year 2015 -> 2016
# ---------------------------------------------------------------------
$wsPax = $wb.Worksheets.Item("BMRESP-passenger")
$wsPax.Range("B1").Value = 2016
[void]$wsPax.Range("A2").Select()

$wsFrt = $wb.Worksheets.Item("BMRESP-freight")
$wsFrt.Range("B1").Value = 2016
[void]$wsFrt.Range("A1").Select()

# ---------------------------------------------------------------------
# 3) Active tab moves from "About" to "BMRESP-freight"
# ---------------------------------------------------------------------
[void]$wsFrt.Activate()
